# EPEX Spot prices workbook — add the "11-nov" day column.
#
# The sheet has one column per day; a new day ("11-nov") needs to be
# inserted right before the "01-oct." column (DN), pushing every
# subsequent day one column to the right (DN..ER -> DO..ES).
# The new column's header gets the date label and its data rows (no
# observations yet for that day) are filled with "-" placeholders,
# matching the convention used for every other missing-data cell in
# this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column at DN, shifting DN:ER (and their contents) one
# column to the right, to DO:ES. Excel/COM semantics: inserting the
# entire column shifts cells right and extends the used range.
$ws.Range("DN1").EntireColumn.Insert()

# New column header (row 1) for the inserted day.
$ws.Range("DN1").Value = "11-nov"

# New column has no price data yet for any hour row -> "-" placeholder,
# same as the rest of the sheet's missing-value cells.
$ws.Range("DN2:DN25").Value = "-"
